$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column A entirely; this shifts columns B:F left to A:E,
# dropping the old column A values (4, 18) and their style.
$ws.Range("A1:A1").EntireColumn.Delete()
